# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" header-name suffixes to the respective
# input-file-derived suffixes "_FV2404" / "_FV2410", turns the used range
# into a native Excel Table ("Table1"), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the column headers in row 1 -------------------------------
$oldToNew = @{
    "Segmentname_old"          = "Segmentname_FV2404"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2404"
    "Segment_old"              = "Segment_FV2404"
    "Datenelement_old"         = "Datenelement_FV2404"
    "Segment ID_old"           = "Segment ID_FV2404"
    "Code_old"                 = "Code_FV2404"
    "Qualifier_old"            = "Qualifier_FV2404"
    "Beschreibung_old"         = "Beschreibung_FV2404"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2404"
    "Bedingung_old"            = "Bedingung_FV2404"
    "Segmentname_new"          = "Segmentname_FV2410"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2410"
    "Segment_new"              = "Segment_FV2410"
    "Datenelement_new"         = "Datenelement_FV2410"
    "Segment ID_new"           = "Segment ID_FV2410"
    "Code_new"                 = "Code_FV2410"
    "Qualifier_new"            = "Qualifier_FV2410"
    "Beschreibung_new"         = "Beschreibung_FV2410"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2410"
    "Bedingung_new"            = "Bedingung_FV2410"
}

$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    if ($oldToNew.ContainsKey($cur)) {
        $cell.Value = $oldToNew[$cur]
    }
}

# --- 2. Turn the data range into a native Excel table ---------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
